$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Choices" column (D) values for rows 2-10 with the new,
# longer option text (replacing the old single-letter placeholder choices).
$ws.Range("D2").Value = "The Office,   TopGear, Sherlock"
$ws.Range("D3").Value = "The Office, TopGear"
$ws.Range("D4").Value = "Sherlock,  TopGear"
$ws.Range("D5").Value = "TopGear"
$ws.Range("D6").Value = "The Grand Tour"
$ws.Range("D7").Value = "E"
$ws.Range("D8").Value = "The Office, Sherlock, TopGear, The Grand Tour"
$ws.Range("D9").Value = "E"
$ws.Range("D10").Value = "Sherlock,  TopGear"

# Widen column D to fit the new, longer choice text.
$ws.Columns.Item(4).ColumnWidth = 37

# Move the active selection as recorded by the editor.
$ws.Range("D23").Select()
